# Update standard deviation measures for the Diastolic Index row (row 4) on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G4").Value = 101
$ws.Range("H4").Value = 1140
$ws.Range("I4").Value = 1022
$ws.Range("J4").Value = 1104
